# Adds two new MIGS.eu.soil.4.0 fields ("culture_collection" and
# "specimen_voucher") as new header columns on row 15, inserted in their
# alphabetically-correct position among the env-package-specific fields,
# each with its DDBJ attribute-description cell comment. Every column
# that sits at/after an insertion point (and its header comment) shifts
# one place to the right, exactly like a normal Excel "Insert Column".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRow = 15
$script:lastCol = 89   # CK15 = last populated header column before this edit

function Find-HeaderColumn($ws, $headerRow, $lastCol, $text) {
    for ($c = 1; $c -le $lastCol; $c++) {
        if ($ws.Cells.Item($headerRow, $c).Value2 -eq $text) {
            return $c
        }
    }
    return -1
}

function Insert-FieldBefore($ws, $headerRow, $beforeFieldName, $newFieldName, $newCommentText) {
    # Column currently holding the field that the new one must precede.
    $insertCol = Find-HeaderColumn $ws $headerRow $script:lastCol $beforeFieldName
    if ($insertCol -eq -1) {
        throw "Could not find header '$beforeFieldName'"
    }

    # Shift existing comments one column to the right, from the last
    # column back down to the insertion point, so nothing is clobbered.
    for ($c = $script:lastCol; $c -ge $insertCol; $c--) {
        $srcCell = $ws.Cells.Item($headerRow, $c)
        $dstCell = $ws.Cells.Item($headerRow, $c + 1)
        $srcComment = $srcCell.Comment
        if ($srcComment -ne $null) {
            $txt = $srcComment.Text()
            $srcComment.Delete()
            $dstCell.AddComment($txt) | Out-Null
        }
    }

    # Shift the header cell values/styles themselves via a real column
    # insert (carries formatting along, like Excel's native behaviour).
    $ws.Columns($insertCol).Insert()

    # Populate the newly-opened column.
    $newCell = $ws.Cells.Item($headerRow, $insertCol)
    $newCell.Value2 = $newFieldName
    $newCell.AddComment($newCommentText) | Out-Null

    $script:lastCol = $script:lastCol + 1
}

$newCommentCulture = @'
Name of source institute and unique culture identifier. See the description for the proper format and list of allowed institutes, http://www.insdc.org/controlled-vocabulary-culturecollection-qualifier
'@

$newCommentVoucher = @'
Identifier for the physical specimen. Use format: "[<institution-code>:[<collection-code>:]]<specimen_id>", eg, "UAM:Mamm:52179". Intended as a reference to the physical specimen that remains after it was analyzed. If the specimen was destroyed in the process of analysis, electronic images (e-vouchers) are an adequate substitute for a physical voucher specimen. Ideally the specimens will be deposited in a curated museum, herbarium, or frozen tissue collection, but often they will remain in a personal or laboratory collection for some time before they are deposited in a curated collection. There are three forms of specimen_voucher qualifiers. If the text of the qualifier includes one or more colons it is a 'structured voucher'. Structured vouchers include institution-codes (and optional collection-codes) taken from a controlled vocabulary maintained by the INSDC that denotes the museum or herbarium collection where the specimen resides, please visit the INSDC website, http://www.insdc.org/controlled-vocabulary-specimenvoucher-qualifier
'@

Insert-FieldBefore $ws $headerRow "cur_land_use" "culture_collection" $newCommentCulture
Insert-FieldBefore $ws $headerRow "store_cond" "specimen_voucher" $newCommentVoucher
